$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells are treated as plain text so numeric-looking
# strings (e.g. "9.50", "581.10") keep their exact formatting/trailing zeros
# instead of being auto-converted to numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.940.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.51%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.570.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.52%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.01'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.37'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.559.58'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.40%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.23%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.219'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +16.92%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.62%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.63'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.46%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.12%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.50%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.134.53'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.81%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.940.90'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.67%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.30'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.35%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.562.01'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.15%  '

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '581.10'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +7.96%  '

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.41'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.62%  '

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'TRON'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.121'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.82%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.01%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.71'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -13.89%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.04'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.58'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.72%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.82'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.30'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.77%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.95'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.49%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.52'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.93%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.27'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.46%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.31'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.29%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.48%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.82'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.36%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.39'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '551.24'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.12%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.417'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.01%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0804'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.52%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.78'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.11%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.05%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +9.95%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.558.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +11.46%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.43'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.75%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0448'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.39%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.21%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.94'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.98%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.37'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.27%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.46%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.02%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.10%  '
